$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.885.45'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.544.37'
$ws.Range("E3").Value = '  -1.17%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.08'
$ws.Range("E5").Value = '  +0.08%  '
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '21.26'
$ws.Range("E9").Value = '  -2.21%  '
$ws.Range("E10").Value = '  -0.22%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0858'
$ws.Range("E11").Value = '  -0.76%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.764.52'
$ws.Range("E12").Value = '  -1.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.543.01'
$ws.Range("E13").Value = '  -1.28%  '
$ws.Range("E14").Value = '  -0.81%  '
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '26.874.63'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.36'
$ws.Range("E17").Value = '  +0.13%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '213.81'
$ws.Range("E18").Value = '  -0.40%  '
$ws.Range("E19").Value = '  +0.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.17'
$ws.Range("E20").Value = '  -2.58%  '
$ws.Range("E21").Value = '  +0.27%  '
$ws.Range("E22").Value = '  -2.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.15'
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("E24").Value = '  -3.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.23'
$ws.Range("E25").Value = '  -1.08%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.59'
$ws.Range("E26").Value = '  -2.14%  '
$ws.Range("E27").Value = '  -0.69%  '
$ws.Range("E28").Value = '  +0.29%  '
$ws.Range("E29").Value = '  +0.08%  '
$ws.Range("E30").Value = '  -1.05%  '
$ws.Range("E31").Value = '  -1.02%  '
$ws.Range("E32").Value = '  +1.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.357.10'
$ws.Range("E33").Value = '  -3.20%  '
$ws.Range("E34").Value = '  +0.54%  '
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.959'
$ws.Range("E36").Value = '  +4.09%  '
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.517'
$ws.Range("E39").Value = '  -1.45%  '
$ws.Range("E40").Value = '  -1.05%  '
$ws.Range("E41").Value = '  +0.26%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.56'
$ws.Range("E42").Value = '  +3.26%  '
$ws.Range("E43").Value = '  -0.77%  '
$ws.Range("E44").Value = '  +1.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.27'
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.73'
$ws.Range("E46").Value = '  -1.56%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.678.62'
$ws.Range("E47").Value = '  -1.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.63'
$ws.Range("E48").Value = '  -0.75%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0510'
$ws.Range("E49").Value = '  +1.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₇0979'
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("E51").Value = '  -0.12%  '
